# Update crypto price/volume figures (Price column D, Volume(1h) column E)
# for rows 2-51 on the active sheet, per the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel
# reinterpret numeric-looking strings (e.g. "1.00", "0.491") as numbers.
# We briefly force a text NumberFormat, assign the value, then restore the
# cell to the default "Normal" style so no stray formatting is introduced.
function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Price column (D) updates ---
$ws.Range("D2").Value = "65.031.17"
$ws.Range("D3").Value = "3.524.65"
Set-TextCell "D5" "592.79"
Set-TextCell "D6" "133.76"
$ws.Range("D7").Value = "3.523.63"
Set-TextCell "D9" "0.491"
Set-TextCell "D10" "0.124"
Set-TextCell "D12" "0.386"
$ws.Range("D13").Value = "4.126.49"
Set-TextCell "D14" "27.68"
$ws.Range("D17").Value = "3.525.64"
$ws.Range("D18").Value = "65.019.40"
Set-TextCell "D19" "10.16"
Set-TextCell "D20" "14.42"
Set-TextCell "D21" "5.69"
Set-TextCell "D22" "392.14"
Set-TextCell "D23" "0.581"
Set-TextCell "D24" "74.95"
$ws.Range("D25").Value = "3.669.39"
Set-TextCell "D27" "0.0000112"
Set-TextCell "D28" "7.75"
Set-TextCell "D32" "8.33"
$ws.Range("D33").Value = "3.532.73"
Set-TextCell "D37" "5.31"
Set-TextCell "D39" "6.96"
Set-TextCell "D40" "168.28"
Set-TextCell "D41" "0.0811"
Set-TextCell "D42" "0.822"
Set-TextCell "D44" "25.81"
Set-TextCell "D45" "42.97"
Set-TextCell "D46" "1.00"
Set-TextCell "D47" "4.44"
Set-TextCell "D48" "1.66"
Set-TextCell "D49" "6.90"
$ws.Range("D50").Value = "2.426.44"
Set-TextCell "D51" "0.909"

# --- Volume(1h) column (E) updates ---
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -3.82%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  +10.70%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  +6.02%  "
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("E43").Value = "  +5.61%  "
$ws.Range("E44").Value = "  -4.31%  "
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").Value = "  +5.72%  "

